# Refresh the crypto price/volume table (and fix the row 50/51 ordering
# for Cosmos / InjectiveProtocol) per the Thu Jul 11 09:51:26 UTC 2024
# GitHub Actions data update.
#
# Note: a couple of "Price" column values (D16, D36, D46, ...) look like
# plain numbers once refreshed (e.g. "4.80", "1.00"); the sheet stores
# them as text, so a leading apostrophe is used to keep Excel from
# re-interpreting them as numeric values and dropping trailing zeros /
# switching to scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.161.54"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "3.132.53"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'529.19"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'142.58"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.135.52"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "3.673.41"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "'25.71"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "'0.0000166"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "58.200.51"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "3.127.14"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "'12.83"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'343.13"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D24").Value = "'0.514"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").Value = "'67.63"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "0.0₃0932"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'7.37"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Value = "'6.41"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "'21.16"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'158.63"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "'4.80"
$ws.Range("E36").Value = "  +4.26%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").Value = "'26.33"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("E40").Value = "  +11.53%  "
$ws.Range("D41").Value = "'0.0667"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'0.701"
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").Value = "3.176.09"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "2.263.07"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  +4.67%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'20.69"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.11"
$ws.Range("E51").Value = "  +1.84%  "
